# Update capital structure database values for rows 2 and 3 (columns G:AQ)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "G"  = -0.002225819506272764
    "H"  = -0.002225819506272764
    "I"  = -0.01070416835289357
    "J"  = -0.01070416835289357
    "K"  = -9.01
    "L"  = -0.01823148522865237
    "M"  = 0.08
    "N"  = 0.008519701810436634
    "O"  = -0.008879023307436182
    "S"  = 0.08
    "U"  = 13.1
    "V"  = 1.395101171458999
    "W"  = -0.3618473895582329
    "X"  = 0.2555988101718658
    "Y"  = -0.6174461997300987
    "Z"  = 4.933120383310042
    "AA" = -0.05280495108804153
    "AB" = 0.05236499156289162
    "AC" = -0.1051699426509332
    "AD" = 62.3
    "AF" = 62.3
    "AG" = 49.2
    "AH" = 0.8690193890361277
    "AI" = 0.7976952624839949
    "AJ" = 0.8397337429595494
    "AK" = 0.7569230769230768
    "AL" = 1.06
    "AM" = 1.06
    "AN" = -14.52214452214452
    "AO" = -4.990566037735849
    "AP" = -11.46853146853147
    "AQ" = -4.990566037735849
}

foreach ($col in $values.Keys) {
    $val = $values[$col]
    $ws.Range("${col}2").Value = $val
    $ws.Range("${col}3").Value = $val
}
